$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 87
$ws.Cells.Item($r, 1).Value = 8
$ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($r, 3).Value = "Coquimbo"
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 4).Value = (Get-Date -Year 2022 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item($r, 5).Value = 4
$ws.Cells.Item($r, 6).Value = 100112030
$ws.Cells.Item($r, 7).Value = "Poroto granado"
$ws.Cells.Item($r, 8).Value = "Sin especificar"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 400
$ws.Cells.Item($r, 11).Value = 27000
$ws.Cells.Item($r, 12).Value = 28000
$ws.Cells.Item($r, 13).Value = 27500
$ws.Cells.Item($r, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item($r, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($r, 16).Value = 1100
$ws.Cells.Item($r, 17).Value = 25
$ws.Cells.Item($r, 18).Value = "Hortaliza"
